$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 126
$ws.Range("A3").Value = 96
$ws.Range("A4").Value = 92
$ws.Range("A5").Value = 90
$ws.Range("A6").Value = 89
$ws.Range("A7").Value = 83
$ws.Range("A8").Value = 78
$ws.Range("A9").Value = 77
$ws.Range("A10").Value = 64
$ws.Range("A11").Value = 23
